$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "year" column (A) held true date values (2016-01-01 / 2016-05-01).
# The edit replaces those with the plain text "2016" in both data rows.
$rng = $ws.Range("A2:A3")
$rng.NumberFormat = "@"
$rng.Value = "2016"
$rng.ClearFormats()
